# feat: ajustes del calculo de cantidad_real
# Update rango_max (column J) so that it matches cantidad_proyectada (column G)
# for every data row of Tabla1 (rows 2-31 on sheet "data").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("data")

$lastRow = 31

for ($r = 2; $r -le $lastRow; $r++) {
    $g = $ws.Cells.Item($r, 7).Value2   # column G = cantidad_proyectada
    $ws.Cells.Item($r, 10).Value2 = $g  # column J = rango_max
}

# Move the active selection to J8, matching the saved view state.
$ws.Range("J8").Select()
